$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$data = @(
    @("SingleUseId28", "Iceland_45", "Left",   "LTR", "WiFi"),
    @("SingleUseId29", "Default",    "Left",   "LTR", "SSID: IMR_TechDemo"),
    @("SingleUseId30", "Default",    "Left",   "LTR", "PASS: Connext123"),
    @("SingleUseId31", "Default",    "Left",   "LTR", "IPV4: 192.168.4.1"),
    @("SingleUseId32", "Default",    "Left",   "LTR", "UPDATE"),
    @("SingleUseId33", "Default",    "Left",   "LTR", "1"),
    @("SingleUseId34", "Default",    "Center", "LTR", "2"),
    @("SingleUseId35", "Default",    "Left",   "LTR", "CONNECTION:")
)

$row = 25
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $row = $row + 1
}
